# A new weekly price-report row needs to be inserted right before the
# existing row 309 (Cilantro / Femacal de La Calera, date 2022-07-08 /
# serial 44754). Inserting there shifts every subsequent record down by
# one row (old 309 -> 310, old 310 -> 311, ... old 384 -> 385), which is
# exactly the pattern shown in the diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 309, pushing rows 309:384 down to 310:385.
$ws.Rows.Item(309).Insert()

# Populate the newly inserted row 309 with the new record's data.
$ws.Cells.Item(309, 1).Value = 3
$ws.Cells.Item(309, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(309, 3).Value = "Coquimbo"
$ws.Cells.Item(309, 4).Value = 44754
$ws.Cells.Item(309, 5).Value = 5
$ws.Cells.Item(309, 6).Value = 100112040
$ws.Cells.Item(309, 7).Value = "Cilantro"
$ws.Cells.Item(309, 8).Value = "Sin especificar"
$ws.Cells.Item(309, 9).Value = "Primera"
$ws.Cells.Item(309, 10).Value = 250
$ws.Cells.Item(309, 11).Value = 3300
$ws.Cells.Item(309, 12).Value = 3500
$ws.Cells.Item(309, 13).Value = 3404
$ws.Cells.Item(309, 14).Value = "`$/docena de atados (3 kilos)"
$ws.Cells.Item(309, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(309, 16).Value = 1135
$ws.Cells.Item(309, 17).Value = 3
$ws.Cells.Item(309, 18).Value = "Hortaliza"
